$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.344528675079346
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 3.873332262039185
$ws.Range("D1").Value = 3.134007930755615
$ws.Range("E1").Value = 1.745621800422668
